# Auto-generated edit script applying numeric updates to Sheets/Ultima_Profits.xlsx
# per the commit diff (scheduled runner data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 2153.3572
$ws.Range("J19").Value = 1715.6
$ws.Range("L19").Value = 1715.6
$ws.Range("N19").Value = -2065.6
$ws.Range("H113").Value = 2327.52
$ws.Range("I113").Value = 2410.4443
$ws.Range("J113").Value = 2114.2856
$ws.Range("K113").Value = 2410.4443
$ws.Range("L113").Value = 2114.2856
$ws.Range("M113").Value = 843.5556999999999
$ws.Range("N113").Value = -8622.285599999999
$ws.Range("H129").Value = 4707.2334
$ws.Range("I129").Value = 371.16666
$ws.Range("J129").Value = 5791.25
$ws.Range("K129").Value = 1113.49998
$ws.Range("L129").Value = 17373.75
$ws.Range("M129").Value = 3886.50002
$ws.Range("N129").Value = -27373.75
$ws.Range("H137").Value = 7408266
$ws.Range("I137").Value = 841.25
$ws.Range("J137").Value = 28572336
$ws.Range("K137").Value = 2523.75
$ws.Range("L137").Value = 85717008
$ws.Range("M137").Value = 26.25
$ws.Range("N137").Value = -85722108
$ws.Range("H138").Value = 2326.7036
$ws.Range("I138").Value = 1216.3334
$ws.Range("J138").Value = 3033.303
$ws.Range("K138").Value = 3649.0002
$ws.Range("L138").Value = 9099.909
$ws.Range("M138").Value = 1490.9998
$ws.Range("N138").Value = -19379.909
$ws.Range("H139").Value = 49800
$ws.Range("J139").Value = 49800
$ws.Range("L139").Value = 49800
$ws.Range("N139").Value = -60080

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H24").Value = 40653.57
$ws.Range("J24").Value = 40653.57
$ws.Range("L24").Value = 40653.57
$ws.Range("N24").Value = -41401.57
$ws.Range("H100").Value = 40653.57
$ws.Range("J100").Value = 40653.57
$ws.Range("L100").Value = 40653.57
$ws.Range("N100").Value = -42817.57
$ws.Range("H139").Value = 77737.5
$ws.Range("J139").Value = 77737.5
$ws.Range("L139").Value = 77737.5
$ws.Range("N139").Value = -88017.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 27412
$ws.Range("J81").Value = 27412
$ws.Range("L81").Value = 27412
$ws.Range("N81").Value = -29534
$ws.Range("H84").Value = 27412
$ws.Range("J84").Value = 27412
$ws.Range("L84").Value = 82236
$ws.Range("N84").Value = -92844
$ws.Range("H134").Value = 3458.64
$ws.Range("I134").Value = 1670.6666
$ws.Range("J134").Value = 6140.6
$ws.Range("K134").Value = 5011.9998
$ws.Range("L134").Value = 18421.8
$ws.Range("M134").Value = -2476.9998
$ws.Range("N134").Value = -23491.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6176127
$ws.Range("I31").Value = 3336.673
$ws.Range("J31").Value = 166668670
$ws.Range("K31").Value = 3336.673
$ws.Range("L31").Value = 166668670
$ws.Range("M31").Value = -3041.673
$ws.Range("N31").Value = -166669260
$ws.Range("H34").Value = 6176127
$ws.Range("I34").Value = 3336.673
$ws.Range("J34").Value = 166668670
$ws.Range("K34").Value = 3336.673
$ws.Range("L34").Value = 166668670
$ws.Range("M34").Value = -3134.673
$ws.Range("N34").Value = -166669074
$ws.Range("H99").Value = 1499.75
$ws.Range("I99").Value = 1383
$ws.Range("J99").Value = 1850
$ws.Range("K99").Value = 1383
$ws.Range("L99").Value = 1850
$ws.Range("M99").Value = 115
$ws.Range("N99").Value = -4846
$ws.Range("H107").Value = 675
$ws.Range("I107").Value = 482.25
$ws.Range("J107").Value = 1960
$ws.Range("K107").Value = 482.25
$ws.Range("L107").Value = 1960
$ws.Range("M107").Value = 1437.75
$ws.Range("N107").Value = -5800
$ws.Range("H126").Value = 1499.75
$ws.Range("I126").Value = 1383
$ws.Range("J126").Value = 1850
$ws.Range("K126").Value = 4149
$ws.Range("L126").Value = 5550
$ws.Range("M126").Value = -1679
$ws.Range("N126").Value = -10490
$ws.Range("H132").Value = 35717356
$ws.Range("I132").Value = 55557460
$ws.Range("J132").Value = 5162.4
$ws.Range("K132").Value = 166672380
$ws.Range("L132").Value = 15487.2
$ws.Range("M132").Value = -166669850
$ws.Range("N132").Value = -20547.2
$ws.Range("H140").Value = 42339.89
$ws.Range("J140").Value = 42339.89
$ws.Range("L140").Value = 42339.89
$ws.Range("N140").Value = -52699.89

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H80").Value = 2750
$ws.Range("I80").Value = 2750
$ws.Range("J80").Value = 2750
$ws.Range("K80").Value = 8250
$ws.Range("L80").Value = 8250
$ws.Range("M80").Value = -7314
$ws.Range("N80").Value = -10122
$ws.Range("H83").Value = 2750
$ws.Range("I83").Value = 2750
$ws.Range("J83").Value = 2750
$ws.Range("K83").Value = 24750
$ws.Range("L83").Value = 24750
$ws.Range("M83").Value = -20070
$ws.Range("N83").Value = -34110
$ws.Range("H105").Value = 2000
$ws.Range("J105").Value = 2000
$ws.Range("L105").Value = 6000
$ws.Range("N105").Value = -11242
$ws.Range("H132").Value = 1197
$ws.Range("I132").Value = 837
$ws.Range("K132").Value = 7533
$ws.Range("M132").Value = -5003  # new cell

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H34").Value = 36200
$ws.Range("J34").Value = 36200
$ws.Range("L34").Value = 36200
$ws.Range("N34").Value = -36736
$ws.Range("H76").Value = 36200
$ws.Range("J76").Value = 36200
$ws.Range("L76").Value = 36200
$ws.Range("N76").Value = -36830
$ws.Range("H79").Value = 36200
$ws.Range("J79").Value = 36200
$ws.Range("L79").Value = 36200
$ws.Range("N79").Value = -38384
$ws.Range("H138").Value = 56176.332
$ws.Range("J138").Value = 56176.332
$ws.Range("L138").Value = 56176.332
$ws.Range("N138").Value = -66456.33199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3818.4565
$ws.Range("I7").Value = 3511.8147
$ws.Range("J7").Value = 4254.2104
$ws.Range("K7").Value = 3511.8147
$ws.Range("L7").Value = 4254.2104
$ws.Range("M7").Value = -3399.8147
$ws.Range("N7").Value = -4478.2104
$ws.Range("H61").Value = 2120
$ws.Range("I61").Value = 1950
$ws.Range("J61").Value = 2233.3333
$ws.Range("K61").Value = 1950
$ws.Range("L61").Value = 2233.3333
$ws.Range("M61").Value = -1748
$ws.Range("N61").Value = -2637.3333
$ws.Range("H113").Value = 2120
$ws.Range("I113").Value = 1950
$ws.Range("J113").Value = 2233.3333
$ws.Range("K113").Value = 1950
$ws.Range("L113").Value = 2233.3333
$ws.Range("M113").Value = 220
$ws.Range("N113").Value = -6573.3333
$ws.Range("H122").Value = 5017.1924
$ws.Range("I122").Value = 4772.35
$ws.Range("K122").Value = 14317.05
$ws.Range("M122").Value = -11867.05
$ws.Range("H126").Value = 3818.4565
$ws.Range("I126").Value = 3511.8147
$ws.Range("J126").Value = 4254.2104
$ws.Range("K126").Value = 10535.4441
$ws.Range("L126").Value = 12762.6312
$ws.Range("M126").Value = -8065.444100000001
$ws.Range("N126").Value = -17702.6312

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 1470.9524
$ws.Range("I113").Value = 1325.6
$ws.Range("J113").Value = 1603.091
$ws.Range("K113").Value = 3976.8
$ws.Range("L113").Value = 4809.272999999999
$ws.Range("M113").Value = -1806.8
$ws.Range("N113").Value = -9149.272999999999
$ws.Range("H126").Value = 4297
$ws.Range("I126").Value = 2138.3635
$ws.Range("J126").Value = 6671.5
$ws.Range("K126").Value = 6415.0905
$ws.Range("L126").Value = 20014.5
$ws.Range("M126").Value = -3945.0905
$ws.Range("N126").Value = -24954.5
$ws.Range("H132").Value = 2357
$ws.Range("I132").Value = 1194.5
$ws.Range("J132").Value = 4100.75
$ws.Range("K132").Value = 3583.5
$ws.Range("L132").Value = 12302.25
$ws.Range("M132").Value = -1053.5
$ws.Range("N132").Value = -17362.25
